$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q7)
$ws.Range("B7").Value = 0.05419696197628752
$ws.Range("C7").Value = 2.10407014602131
$ws.Range("D7").Value = 13.00453704867561
$ws.Range("E7").Value = 3.606180396025081
$ws.Range("F7").Value = 3.654174917201057
$ws.Range("G7").Value = 38

# Row 8 (Q8)
$ws.Range("B8").Value = 0.228655978109646
$ws.Range("C8").Value = 2.04971899904124
$ws.Range("D8").Value = 12.64796182215585
$ws.Range("E8").Value = 3.55639730937867
$ws.Range("F8").Value = 3.597993637801378
$ws.Range("G8").Value = 37

# Row 9 (Q9)
$ws.Range("B9").Value = 0.2908677701452642
$ws.Range("C9").Value = 2.547516861569405
$ws.Range("D9").Value = 20.99274268013081
$ws.Range("E9").Value = 4.581783788016498
$ws.Range("F9").Value = 4.691328912884362
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("B10").Value = -0.9702732014409903
$ws.Range("C10").Value = 2.007183996697524
$ws.Range("D10").Value = 8.368548844988011
$ws.Range("E10").Value = 2.892844421151613
$ws.Range("F10").Value = 2.836555185934091
$ws.Range("G10").Value = 13

# Row 11
$ws.Range("B11").Value = -0.934137815474154
$ws.Range("C11").Value = 1.739680744534949
$ws.Range("D11").Value = 5.346651227107502
$ws.Range("E11").Value = 2.312282687542226
$ws.Range("F11").Value = 2.364856699889202
$ws.Range("G11").Value = 5
